# repull data, push all data, mean calculation
# Update the dSF (column F) values for the rows where the source data was repulled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -11
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = 9
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = 3
$ws.Range("F10").Value = -3
$ws.Range("F11").Value = -4
$ws.Range("F12").Value = -3
$ws.Range("F13").Value = -1
$ws.Range("F14").Value = -6
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = 3
